# Applies the "right and wrong kill handlers" edit to shanRoyale2022Data1.xlsx
$wb = $excel.ActiveWorkbook

# Target sheet is the first sheet (playerDataRound1), which was the active/selected sheet.
$ws = $wb.Worksheets.Item("playerDataRound1")
$ws.Activate()

# Update faction (C4) and dying (D4) values for row 4 (vigonometry / VIGNESHWAR).
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 1

# Move the active selection to C5 (matches saved cursor position in the file).
$ws.Range("C5").Select()

# Best-effort: restore the window geometry recorded by Excel when the file was saved.
$win = $excel.ActiveWindow
$win.Left = -27660
$win.Top = 1140
$win.Width = 21600
$win.Height = 11175
